$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D width (Desktop work: widened description column) ---
$ws.Columns.Item(4).ColumnWidth = 47.86

# --- Row 3 untouched (times already set) ---

# --- Row 4: new research task, entered 11/1/2019 ---
$ws.Rows.Item(4).RowHeight = 98.25

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "date:11/1/2019"
$ws.Range("D4").Value = "Research about technologies which are required for following!`n 1. front-end:Xml,html,css`n 2. back-end:java,php,python,`n 3. server:aws,Apache,Firebase`n 4. database:MongoDB,SQLite,Mysql,Firebase(Real time database) etc.."
$ws.Range("E4").Value = "search about technology required for front-end,bach-end,server,database of application "
$ws.Range("E5").Value = "search abiut the formate of quiz"
$ws.Range("F4").Value = "https://www.quiz-zone.co.uk/runningaquiz/format.html"
$ws.Range("F5").Value = "https://cdl.ucf.edu/support/webcourses/respondus/quiz-formatting-guidelines/"
$ws.Range("D5").Value = "prototype of quiz app in adob xd"

# --- Row 5: second task of the day, with start/end time ---
$ws.Rows.Item(5).RowHeight = 26.25

$ws.Range("B5").NumberFormat = "h:mm"
$ws.Range("B5").Value = 0.39583333333333331
$ws.Range("C5").NumberFormat = "h:mm"
$ws.Range("C5").Value = 0.47916666666666669

# --- Update selection to land on D5, matching the active cell at save time ---
[void]$ws.Range("D5").Select()
